$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-30 Saturday", "2024-12-01 Sunday"),
    @("267×2=534", "616×2=1232"),
    @("212×7=1484", "616×5=3080"),
    @("774×8=6192", "393×4=1572"),
    @("177×7=1239", "306×4=1224"),
    @("763×3=2289", "711×8=5688"),
    @("229×5=1145", "778×8=6224"),
    @("712×6=4272", "595×9=5355"),
    @("408×6=2448", "254×6=1524"),
    @("586×8=4688", "733×5=3665"),
    @("455×3=1365", "513×2=1026"),
    @("241×4=964", "148×6=888"),
    @("788×4=3152", "559×3=1677"),
    @("205×5=1025", "892×3=2676"),
    @("619×6=3714", "316×7=2212"),
    @("308×2=616", "648×8=5184"),
    @("748×3=2244", "728×7=5096"),
    @("933×7=6531", "907×4=3628"),
    @("150×4=600", "108×5=540"),
    @("824×5=4120", "263×9=2367"),
    @("370×5=1850", "931×3=2793"),
    @("300×3=900", "574×2=1148"),
    @("265×4=1060", "191×7=1337"),
    @("223×8=1784", "317×2=634"),
    @("533×2=1066", "292×6=1752"),
    @("485×9=4365", "165×9=1485")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
